$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 3552926.13
$ws.Range("C7").Value = -20.03457886533852
$ws.Range("D7").Value = 3107
$ws.Range("E7").Value = 3107
$ws.Range("F7").Value = 1143.523054393305
$ws.Range("G7").Value = 21.89128886184639
